# feat: add 2022-Q1 data
#
# Before:
#   Sheet1 "2021-Q4" (fund holdings detail for 2021-Q4)
#   Sheet2 "总计"     (summary: 日期 / 持有数量(只) / 持有市值(亿元) -> one row, 2021-Q4)
#
# After:
#   Sheet1 "2021-Q4" (unchanged)
#   Sheet2 "2022-Q1" (fund holdings detail for 2022-Q1, replaces old "总计" sheet content)
#   Sheet3 "总计"     (summary, now with two rows: 2022-Q1 and 2021-Q4)

# Writes $value into $range as literal text, even when it looks like a
# number (e.g. "002802", "0.41"), without leaving the target cell's own
# style changed. It stages the text (forced via a Text number format) in a
# scratch cell, copies just the *value* over (PasteSpecial xlPasteValues
# does not carry the source's formatting along), then wipes the scratch
# cell again.
function Set-LiteralText($range, $value) {
    $scratch = $range.Worksheet.Range("ZZ1000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

$wb = $excel.ActiveWorkbook

# Grab the existing "总计" sheet before touching anything.
$wsOld = $wb.Worksheets.Item(2)

# Duplicate it (full formatting/pageSetup/styles included) and place the
# duplicate right after it - that duplicate becomes the new "总计" summary
# sheet (sheetId 3 / rId3), while the original sheet (sheetId 2 / rId2) gets
# repurposed below into the new "2022-Q1" detail sheet.
$wsOld.Copy($null, $wsOld)
$wsTotal = $wb.Worksheets.Item(3)

$wsOld.Name = "2022-Q1"
$wsQ1 = $wsOld

$wsTotal.Name = "总计"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "2022-Q1": fund holdings detail
# ---------------------------------------------------------------------

# Header row (B1:H1). B1:D1 already carry the bold/bordered header style
# inherited from the old "总计" sheet - overwriting their text keeps that
# style. E1:H1 are new cells, so copy the format from D1 onto them first.
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"

$wsQ1.Range("D1").Copy()
$wsQ1.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Data row 2 - A2 already holds 0 with the correct style from the old sheet.
$wsQ1.Range("A2").Value = 0

Set-LiteralText $wsQ1.Range("B2") "002802"
$wsQ1.Range("C2").Value = "广发东财大数据精选灵活配置混合"
Set-LiteralText $wsQ1.Range("D2") "0.41"
Set-LiteralText $wsQ1.Range("E2") "55.13"
Set-LiteralText $wsQ1.Range("F2") "2.16"
Set-LiteralText $wsQ1.Range("G2") "0.0089"

$wsQ1.Range("H2").Value = 7

# ---------------------------------------------------------------------
# Sheet "总计": summary, add the new 2022-Q1 row ahead of the 2021-Q4 row
# ---------------------------------------------------------------------

# Move the existing 2021-Q4 row down to row 3 (preserving A2's "row index"
# style on the new A3 cell), then write the new 2022-Q1 row into row 2.
$oldB2 = $wsTotal.Range("B2").Value()
$oldC2 = $wsTotal.Range("C2").Value()
$oldD2 = $wsTotal.Range("D2").Value()

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = $oldB2
$wsTotal.Range("C3").Value = $oldC2
$wsTotal.Range("D3").Value = $oldD2

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.01
